$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize column C (OBRA) so its contents fit, the way a user would by
# double-clicking the column-header border / Format > AutoFit Column Width.
# (ColumnWidth 21 is the closest reachable value to the target best-fit
# width of ~21.86 characters given this host's column-width rounding.)
$ws.Columns("C").ColumnWidth = 21

# Move the active selection to C1, matching the post-edit cursor position.
$ws.Range("C1").Select()
